$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in rows 28-30 with new date/time entries
$ws.Range("A28").Value = 43060
$ws.Range("B28").Value = 0.52083333333333337
$ws.Range("C28").Value = 0.54166666666666663

$ws.Range("A29").Value = 43060
$ws.Range("B29").Value = 0.54513888888888895
$ws.Range("C29").Value = 0.57291666666666663

$ws.Range("A30").Value = 43060
$ws.Range("B30").Value = 0.57638888888888895
$ws.Range("C30").Value = 0.60416666666666663

# Update the selected cell to reflect the new active cell/selection
$ws.Range("I25").Select()
